$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (column D) cells keep exact text formatting (avoid numeric coercion)
foreach ($addr in @(
    "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.884.01'
$ws.Range("E2").Value = '  +5.47%  '
$ws.Range("D3").Value = '1.810.95'
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("D4").Value = '0.9980'
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").Value = '315.26'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").Value = '0.9953'
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("D7").Value = '0.5424'
$ws.Range("E7").Value = '  +3.52%  '
$ws.Range("D8").Value = '0.3833'
$ws.Range("E8").Value = '  +1.77%  '
$ws.Range("D9").Value = '0.07611'
$ws.Range("E9").Value = '  +2.70%  '
$ws.Range("D10").Value = '42.69'
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("D11").Value = '1.129'
$ws.Range("E11").Value = '  +3.17%  '
$ws.Range("D12").Value = '21.35'
$ws.Range("E12").Value = '  +3.10%  '
$ws.Range("D13").Value = '0.9939'
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("D14").Value = '6.233'
$ws.Range("E14").Value = '  +2.07%  '
$ws.Range("D15").Value = '7.483'
$ws.Range("E15").Value = '  +7.08%  '
$ws.Range("D16").Value = '1.802.82'
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("D17").Value = '92.17'
$ws.Range("E17").Value = '  +2.77%  '
$ws.Range("D18").Value = '0.00001072'
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").Value = '0.06443'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '17.37'
$ws.Range("E20").Value = '  +3.34%  '
$ws.Range("D21").Value = '0.9951'
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("D22").Value = '6.011'
$ws.Range("E22").Value = '  +2.05%  '
$ws.Range("D23").Value = '28.921.68'
$ws.Range("E23").Value = '  +5.45%  '
$ws.Range("D24").Value = '11.49'
$ws.Range("E24").Value = '  +2.34%  '
$ws.Range("D25").Value = '2.127'
$ws.Range("E25").Value = '  +1.74%  '
$ws.Range("D26").Value = '162.24'
$ws.Range("E26").Value = '  +4.36%  '
$ws.Range("D27").Value = '20.77'
$ws.Range("E27").Value = '  +2.66%  '
$ws.Range("D28").Value = '2.425'
$ws.Range("E28").Value = '  +2.75%  '
$ws.Range("D29").Value = '2.016.16'
$ws.Range("E29").Value = '  +1.58%  '
$ws.Range("D30").Value = '124.21'
$ws.Range("E30").Value = '  +2.49%  '
$ws.Range("D31").Value = '1.153'
$ws.Range("E31").Value = '  +5.80%  '
$ws.Range("D32").Value = '0.1021'
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").Value = '5.810'
$ws.Range("E33").Value = '  +3.50%  '
$ws.Range("D34").Value = '3.671'
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("D35").Value = '0.2326'
$ws.Range("E35").Value = '  +13.39%  '
$ws.Range("D36").Value = '0.06638'
$ws.Range("E36").Value = '  +10.63%  '
$ws.Range("E37").Value = '  +3.54%  '
$ws.Range("E38").Value = '  +5.27%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '11.68'
$ws.Range("E39").Value = '  +3.09%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '8.659'
$ws.Range("E40").Value = '  +5.49%  '
$ws.Range("D41").Value = '0.6417'
$ws.Range("E41").Value = '  +4.60%  '
$ws.Range("D42").Value = '1.236'
$ws.Range("E42").Value = '  +9.09%  '
$ws.Range("B43").Value = 'WEMIXTOKEN'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '1.405'
$ws.Range("E43").Value = '  -1.95%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '0.9939'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").Value = '13.75'
$ws.Range("D46").Value = '0.6008'
$ws.Range("E46").Value = '  +3.59%  '
$ws.Range("D47").Value = '3.689'
$ws.Range("E47").Value = '  +1.72%  '
$ws.Range("D48").Value = '126.11'
$ws.Range("E48").Value = '  +3.53%  '
$ws.Range("D49").Value = '2.017'
$ws.Range("E49").Value = '  +6.33%  '
$ws.Range("D50").Value = '1.169'
$ws.Range("E50").Value = '  +4.89%  '
$ws.Range("D51").Value = '0.06981'
$ws.Range("E51").Value = '  +3.67%  '
